$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(44441,44442,44443,44444,44445,44446,44447,44448)
$col_b = @(5,6,5,4,0,2,1,8)
$col_c = @(39,37,38,30,26,23,23,26)
$col_d = @(152.6478531449372,144.8197581118635,148.7338056284003,117.4214254961055,101.7652354299581,90.02309288034756,90.02309288034756,101.7652354299581)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 367 + $i

    # Copy the formatting (including the date number format + border/font style)
    # from the last existing data row (366) down onto column A of the new row.
    $ws.Cells.Item(366, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $col_b[$i]
    $ws.Cells.Item($row, 3).Value = $col_c[$i]
    $ws.Cells.Item($row, 4).Value = $col_d[$i]
}

$excel.CutCopyMode = $false
